$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for specific rows per repulled data / mean calc fix
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = -8
$ws.Range("F8").Value = -6
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -7
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = -3
$ws.Range("F15").Value = -5
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("F19").Value = -2
